# Updated cryptos list data (Price and Volume(1h) columns)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) holds numeric-looking text (e.g. "228.01"); force it to stay
# text so Excel does not silently convert it to a floating-point number.
$priceCells = @("D2", "D3", "D5", "D6", "D8", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D19", "D20", "D21", "D25", "D26", "D27", "D35", "D36", "D37", "D40", "D41", "D48", "D49")
foreach ($cellRef in $priceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = '34.523.82'
$ws.Range("E2").Value = '  -0.22%  '
$ws.Range("D3").Value = '1.808.03'
$ws.Range("E3").Value = '  -0.56%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").Value = '228.01'
$ws.Range("E5").Value = '  -0.15%  '
$ws.Range("D6").Value = '0.581'
$ws.Range("E6").Value = '  +3.94%  '
$ws.Range("E7").Value = '  +0.13%  '
$ws.Range("D8").Value = '36.83'
$ws.Range("E8").Value = '  +5.84%  '
$ws.Range("E9").Value = '  -0.59%  '
$ws.Range("D10").Value = '0.0694'
$ws.Range("E10").Value = '  -0.33%  '
$ws.Range("D11").Value = '0.0965'
$ws.Range("E11").Value = '  +1.42%  '
$ws.Range("D12").Value = '2.068.39'
$ws.Range("E12").Value = '  -0.59%  '
$ws.Range("D13").Value = '11.57'
$ws.Range("E13").Value = '  +1.75%  '
$ws.Range("D14").Value = '0.654'
$ws.Range("E14").Value = '  +1.38%  '
$ws.Range("D15").Value = '1.801.39'
$ws.Range("E15").Value = '  -1.12%  '
$ws.Range("D16").Value = '4.49'
$ws.Range("E16").Value = '  +3.33%  '
$ws.Range("D17").Value = '34.491.95'
$ws.Range("E17").Value = '  -0.28%  '
$ws.Range("E18").Value = '  +0.59%  '
$ws.Range("D19").Value = '245.47'
$ws.Range("E19").Value = '  -0.76%  '
$ws.Range("D20").Value = '0.0₃0792'
$ws.Range("E20").Value = '  -1.33%  '
$ws.Range("D21").Value = '11.68'
$ws.Range("E21").Value = '  +1.00%  '
$ws.Range("E22").Value = '  +0.14%  '
$ws.Range("E23").Value = '  -0.43%  '
$ws.Range("E24").Value = '  +5.29%  '
$ws.Range("D25").Value = '172.68'
$ws.Range("E25").Value = '  -0.10%  '
$ws.Range("D26").Value = '7.98'
$ws.Range("E26").Value = '  +7.06%  '
$ws.Range("D27").Value = '16.94'
$ws.Range("E27").Value = '  +0.96%  '
$ws.Range("E28").Value = '  +1.08%  '
$ws.Range("E29").Value = '  +0.11%  '
$ws.Range("E30").Value = '  -0.70%  '
$ws.Range("E31").Value = '  -0.23%  '
$ws.Range("E32").Value = '  +0.05%  '
$ws.Range("E33").Value = '  -0.46%  '
$ws.Range("E34").Value = '  -1.70%  '
$ws.Range("D35").Value = '1.396.35'
$ws.Range("E35").Value = '  -1.63%  '
$ws.Range("D36").Value = '0.672'
$ws.Range("E36").Value = '  -0.69%  '
$ws.Range("D37").Value = '2.45'
$ws.Range("E37").Value = '  -5.29%  '
$ws.Range("E38").Value = '  -0.63%  '
$ws.Range("E39").Value = '  -0.41%  '
$ws.Range("D40").Value = '0.970'
$ws.Range("E40").Value = '  +1.59%  '
$ws.Range("D41").Value = '83.24'
$ws.Range("E41").Value = '  -3.30%  '
$ws.Range("E42").Value = '  -0.62%  '
$ws.Range("E43").Value = '  +0.95%  '
$ws.Range("E44").Value = '  +7.93%  '
$ws.Range("E45").Value = '  -2.28%  '
$ws.Range("E46").Value = '  -1.18%  '
$ws.Range("E47").Value = '  -4.88%  '
$ws.Range("D48").Value = '1.969.53'
$ws.Range("E48").Value = '  -0.61%  '
$ws.Range("D49").Value = '104.59'
$ws.Range("E49").Value = '  -1.12%  '
$ws.Range("E50").Value = '  +0.14%  '
$ws.Range("E51").Value = '  -2.73%  '

# Restore default (Normal) style now that the text values are safely stored, so no
# stray "@" number-format styling is left attached to these cells.
foreach ($cellRef in $priceCells) {
    $ws.Range($cellRef).Style = "Normal"
}
